$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.18"
$ws.Range("D3").Value = "'23.05"
$ws.Range("D4").Value = "'5.402"
$ws.Range("D5").Value = "'0.05985"
$ws.Range("D7").Value = "'6.489"
$ws.Range("D8").Value = "'0.8081"
$ws.Range("D9").Value = "'0.9242"
$ws.Range("D10").Value = "'0.1429"
$ws.Range("D11").Value = "'0.07419"
$ws.Range("D12").Value = "'0.03266"
$ws.Range("D13").Value = "'0.03074"
$ws.Range("D14").Value = "'0.09360"
$ws.Range("D15").Value = "'3.850"
$ws.Range("D16").Value = "'0.001583"
$ws.Range("D17").Value = "'0.04699"
$ws.Range("D18").Value = "'0.01112"
$ws.Range("E18").Value = "17OneONEBestin24h"
$ws.Range("D19").Value = "'0.005867"
$ws.Range("D20").Value = "'0.001267"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("D21").Value = "'0.004879"
$ws.Range("D22").Value = "'0.00006799"
$ws.Range("D23").Value = "'3.572"
$ws.Range("D24").Value = "'2.178"
$ws.Range("D25").Value = "'0.3235"
$ws.Range("D26").Value = "'0.1329"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006369"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.004299"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1076"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008920"
$ws.Range("D45").Value = "'0.00005086"
$ws.Range("D47").Value = "'0.6499"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
